# Adding documentation Test scripts
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("magento_new_users")

# Update test data values
$ws.Range("D2").Value = "mando24052025@gmail.com"
$ws.Range("D3").Value = "mando26052025@gmail.com"
$ws.Range("G3").Value = "Pass - User Created"
$ws.Range("G4").Value = "Pass - First and last Name are not valid"

# Update the active selection on the sheet
$ws.Range("C18:D20").Select()
